$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: Name / email / date (with trailing spaces as in source)
$ws.Range("F4").Value = "Name "
$ws.Range("H4").Value = "email "
$ws.Range("J4").Value = "date "

# Box a thin border around the small F4:J5 table
$ws.Range("F4:J5").Borders.LineStyle = 1

# Leave the cursor parked on G10, matching the saved selection
$ws.Range("G10").Select()
